$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos-list refresh: update Price (D) / Volume(1h) (E) for every
# coin row, and for the handful of rows whose rank order changed, also
# rewrite Coin (B) and Link (C) so the row holds the new coin.
#
# Some new Price values are plain numeric-looking strings (e.g. "537.63",
# "0.110"). Setting .Value directly would make Excel coerce them to real
# numbers (losing significant trailing zeros, e.g. "0.110" -> 0.11), so
# those specific cells are pre-formatted as Text first.

# Row 2
$ws.Range("D2").Value = "61.624.24"
$ws.Range("E2").Value = "  -3.39%  "

# Row 3
$ws.Range("D3").Value = "3.002.68"

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.63"
$ws.Range("E5").Value = "  -0.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.07"
$ws.Range("E6").Value = "  -2.43%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "2.997.38"
$ws.Range("E8").Value = "  -2.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("E10").Value = "  -5.27%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -3.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  -2.73%  "

# Row 13
$ws.Range("E13").Value = "  -2.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.85"
$ws.Range("E14").Value = "  -2.99%  "

# Row 15
$ws.Range("D15").Value = "3.489.97"

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "61.652.89"
$ws.Range("E16").Value = "  -3.37%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.110"
$ws.Range("E17").Value = "  -1.57%  "

# Row 18
$ws.Range("D18").Value = "3.003.19"
$ws.Range("E18").Value = "  -2.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.61"
$ws.Range("E19").Value = "  -1.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.75"
$ws.Range("E20").Value = "  -4.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.23"
$ws.Range("E21").Value = "  -2.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -4.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.92"
$ws.Range("E23").Value = "  -3.72%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.50"
$ws.Range("E24").Value = "  +0.72%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.95"
$ws.Range("E25").Value = "  -2.60%  "

# Row 26
$ws.Range("E26").Value = "  -0.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -1.75%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.73"
$ws.Range("E28").Value = "  -7.29%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  +3.19%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.61"
$ws.Range("E31").Value = "  -2.70%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.88"
$ws.Range("E32").Value = "  -1.74%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.46"
$ws.Range("E33").Value = "  -0.83%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.18"
$ws.Range("E34").Value = "  -2.80%  "

# Row 35
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.28"
$ws.Range("E35").Value = "  -5.78%  "

# Row 36
$ws.Range("E36").Value = "  -3.12%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "454.06"
$ws.Range("E37").Value = "  -8.06%  "

# Row 38
$ws.Range("D38").Value = "3.169.60"
$ws.Range("E38").Value = "  -4.41%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0785"
$ws.Range("E39").Value = "  -2.11%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.120"
$ws.Range("E40").Value = "  +1.51%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0384"
$ws.Range("E41").Value = "  -4.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.09"
$ws.Range("E42").Value = "  -1.12%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("E43").Value = "  -8.45%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.53"
$ws.Range("E44").Value = "  +5.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.243"
$ws.Range("E46").Value = "  -6.73%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").Value = "  -5.83%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.55"
$ws.Range("E48").Value = "  -3.12%  "

# Row 49
$ws.Range("E49").Value = "  -1.58%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0495"
$ws.Range("E50").Value = "  -8.38%  "

# Row 51
$ws.Range("E51").Value = "  +5.86%  "
